$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new row 4 below the table, copying row 3's formatting as a starting point ---
$ws.Range("A3:D3").Copy()
$ws.Range("A4:D4").PasteSpecial(-4122)
$ws.Rows("4").RowHeight = 30

# --- Extend the A column merge down through the new row 4 ---
$ws.Range("A2:A3").UnMerge()
$ws.Range("A2:A4").Merge()

# --- Column C becomes wider (~40.71) while A, B, D stay at 20.7109375 ---
$ws.Columns("C").ColumnWidth = 39.75

# --- Row 2: B2 becomes "Life" using the plain style (like C2/A2), C2 left blank, then merge B2:C2 ---
$ws.Range("A3").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B2").Value = "Life"
$ws.Range("C2").Value = ""
$ws.Range("B2:C2").Merge()

# --- Row 3: B3 becomes "Generation" using the same plain style, C3 blank, then merge B3:C3 ---
$ws.Range("A3").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B3").Value = "Generation"
$ws.Range("C3").Value = ""
$ws.Range("B3:C3").Merge()

# --- Row 4: B4 = "Age" (highlighted style, same as old B2/B3), C4 = "Birth" (plain style) ---
$ws.Range("B4").Value = "Age"
$ws.Range("C4").Value = "Birth"

# --- Extend the D column merge down through the new row 4 (added last to match merge order) ---
$ws.Range("D2:D3").UnMerge()
$ws.Range("D2:D4").Merge()
